{"js": "const body = context.document.body;\n\n// --- Change 1 -------------------------------------------------------\n// \"...b\u1eb1ng c\u00e2u l\u1ec7nh:\" -> \"...b\u1eb1ng c\u00e2u l\u1ec7nh, m\u1eadt kh\u1ea9u l\u00e0 \u201cubuntu\u201d:\"\n// Find the unique phrase right before the trailing colon (excluding the\n// colon itself) and insert the new text right after it, i.e. right\n// before the colon.\nconst results1 = body.search(\"sang cho bob b\u1eb1ng c\u00e2u l\u1ec7nh\", { matchCase: true });\nresults1.load(\"text\");\nawait context.sync();\n\nif (results1.items.length > 0) {\n  results1.items[0].insertText(\", m\u1eadt kh\u1ea9u l\u00e0 \\u201Cubuntu\\u201D\", \"End\");\n  await context.sync();\n}\n\n// --- Change 2 -------------------------------------------------------\n// \"... -t .raw -)\" -> \"... -t .raw -) | head -n 20\"\n// Find the unique phrase that ends the command line and append the new\n// text right after it.\nconst results2 = body.search(\"secret.wav -t .raw -)\", { matchCase: true });\nresults2.load(\"text\");\nawait context.sync();\n\nif (results2.items.length > 0) {\n  results2.items[0].insertText(\" | head -n 20\", \"End\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Change 1 --------------------------------------------------------\n# \"...b\u1eb1ng c\u00e2u l\u1ec7nh:\" -> \"...b\u1eb1ng c\u00e2u l\u1ec7nh, m\u1eadt kh\u1ea9u l\u00e0 \u201cubuntu\u201d:\"\n# Find the unique phrase right before the trailing colon (excluding the\n# colon itself), collapse to its end point and insert the new text\n# right there (i.e. right before the colon).\n$rng1 = $d.Content\n$find1 = $rng1.Find\n$find1.ClearFormatting()\n$find1.Text = \"sang cho bob b\u1eb1ng c\u00e2u l\u1ec7nh\"\n$find1.MatchCase = $true\n$found1 = $find1.Execute()\nif ($found1) {\n    $rng1.Collapse(0)  # wdCollapseEnd\n    $rng1.Text = \", m\u1eadt kh\u1ea9u l\u00e0 \u201cubuntu\u201d\"\n}\n\n# --- Change 2 --------------------------------------------------------\n# \"... -t .raw -)\" -> \"... -t .raw -) | head -n 20\"\n# Find the unique phrase that ends the command line, collapse to its\n# end point and insert the new text right after it.\n$rng2 = $d.Content\n$find2 = $rng2.Find\n$find2.ClearFormatting()\n$find2.Text = \"secret.wav -t .raw -)\"\n$find2.MatchCase = $true\n$found2 = $find2.Execute()\nif ($found2) {\n    $rng2.Collapse(0)  # wdCollapseEnd\n    $rng2.Text = \" | head -n 20\"\n}\n"}
